# Outstandings.xlsx update
#  - Sheet "Purchase 22-23": group 1 (Namrata Rubber) gains a third invoice line,
#    the original first line ("100/23-24") is dropped, the remaining two lines'
#    dates/invoice numbers/amounts shift, and a brand-new group 5 (Pilz India)
#    is appended.
#  - Sheet "Sale 22-23": group 2 (Putzmeister) gains a new invoice line.
#  - The active sheet/selection flips from "Sale 22-23" back to "Purchase 22-23".

$wb = $excel.ActiveWorkbook

$wsPurchase = $wb.Worksheets.Item("Purchase 22-23")
$wsSale     = $wb.Worksheets.Item("Sale 22-23")

# ---------------------------------------------------------------------------
# Sheet "Purchase 22-23"
# ---------------------------------------------------------------------------

# Insert a new row before the existing blank separator (old row 5), shifting
# everything below it down by one. Formulas referencing the shifted rows are
# adjusted automatically by Excel.
$wsPurchase.Rows.Item(5).Insert()

# Row 3 (still row 3): date / invoice no. / outstanding-bill amount change.
$wsPurchase.Range("B3").Value = 45355
$wsPurchase.Range("C3").Value = "114/23-24"
$wsPurchase.Range("E3").Value = 47466

# Row 4 (still row 4): date / invoice no. / bill amount change, and the
# running-total formula that used to live here is removed (moves to row 5).
$wsPurchase.Range("B4").Value = 45363
$wsPurchase.Range("C4").Value = "119/23-24"
$wsPurchase.Range("E4").Value = 10178
$wsPurchase.Range("F4").ClearContents()

# New row 5: third invoice line for the same vendor, with the running total.
$wsPurchase.Range("A4:F4").Copy()
$wsPurchase.Range("A5:F5").PasteSpecial(-4122) | Out-Null
$wsPurchase.Range("B5").Value = 45365
$wsPurchase.Range("C5").Value = "122/23-24"
$wsPurchase.Range("D5").Value = "Namrata Rubber Product Pvt Ltd"
$wsPurchase.Range("E5").Value = 129151
$wsPurchase.Range("F5").Formula = "=E3+E4+E5"

# New group 5 (row 17, two rows below the last existing row 15 - row 16 stays
# blank/unused just like row 14 and row 16 are skipped elsewhere in the sheet).
$wsPurchase.Range("A11:F11").Copy()
$wsPurchase.Range("A17:F17").PasteSpecial(-4122) | Out-Null
$wsPurchase.Range("A17").Value = 5
$wsPurchase.Range("B17").Value = 45364
$wsPurchase.Range("C17").Value = "INV-017490"
$wsPurchase.Range("D17").Value = "Pilz India Pvt Ltd"
$wsPurchase.Range("E17").Value = 192635
$wsPurchase.Range("F17").Formula = "=E17"

$wsPurchase.Range("F6").Select()

# ---------------------------------------------------------------------------
# Sheet "Sale 22-23"
# ---------------------------------------------------------------------------

# Insert a new row before the blank separator (old row 9), shifting everything
# below it down by one; the running-total formula auto-adjusts to the new row.
$wsSale.Rows.Item(9).Insert()

$wsSale.Range("A8:F8").Copy()
$wsSale.Range("A9:F9").PasteSpecial(-4122) | Out-Null
$wsSale.Range("A9").Value = ""
$wsSale.Range("B9").Value = 45370
$wsSale.Range("C9").Value = "b23-24MQ409"
$wsSale.Range("D9").Value = "Putzmeister Concrete Machines Pvt Ltd"
$wsSale.Range("E9").Value = 793578
$wsSale.Range("F9").Formula = "=E5+E6+E7+E8+E9"

# The formula that used to total the group now lives on row 9; row 8 no
# longer carries a value.
$wsSale.Range("F8").ClearContents()

$wsSale.Range("G14").Select()

# ---------------------------------------------------------------------------
# Active sheet / tab selection flips back to "Purchase 22-23"
# ---------------------------------------------------------------------------
$wsPurchase.Activate()
$wsPurchase.Select()
